$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")

# Insert a new row above row 3 ("Arbeitnehmerbeitrag gesetzliche Pflegeversicherung...")
# to hold the new "Juenger als 23 oder geboren vor 1940" boolean field.
$ws1.Rows.Item(3).Insert()

# Populate the new row 3 on Tabelle1 (label only for now; the value is set
# after the shared strings for the lookup sheet so the shared-string table
# keeps the same add-order as the source workbook).
$ws1.Range("A3").Value = "Juenger als 23 oder geboren vor 1940"

# New hidden lookup sheet "Tabelle2" holding the boolean list values.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Tabelle2"
$ws2.Range("A1").Value = "boolean"
$ws2.Range("A2").Value = "ja"
$ws2.Range("A3").Value = "nein"
$ws2.Visible = [Microsoft.Office.Interop.Excel.XlSheetVisibility]::xlSheetHidden

$ws1.Range("B3").Value = "nein"

# Add the dropdown list data validation on B3, sourced from Tabelle2!$A$2:$A$3.
$validation = $ws1.Range("B3").Validation
$validation.Add(
    [Microsoft.Office.Interop.Excel.XlDVType]::xlValidateList,
    [Microsoft.Office.Interop.Excel.XlFormatConditionOperator]::xlBetween,
    [Microsoft.Office.Interop.Excel.XlFormatConditionOperator]::xlBetween,
    "=Tabelle2!`$A`$2:`$A`$3"
)
$validation.IgnoreBlank = $true
$validation.InCellDropdown = $true
$validation.ShowInput = $true
$validation.ShowError = $true

$ws1.Activate()
$ws1.Range("A9").Select()
